$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.921.84"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.846.85"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "699.46"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.00"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "3.845.76"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.24"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "4.493.68"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "3.881.39"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "70.925.22"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.55"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.56"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.44"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").Value = "3.801.61"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("E40").Value = "  +6.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.00"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -9.50%  "
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.72"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.34"
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("E51").Value = "  -4.86%  "
